$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 115, shifting rows 115:216 down to 116:217
$ws.Rows.Item(115).Insert()

# Populate the new row 115 with the new price record
$ws.Cells.Item(115, 1).Value  = 11
$ws.Cells.Item(115, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(115, 3).Value  = "Bíobío"
$ws.Cells.Item(115, 4).Value  = 44589
$ws.Cells.Item(115, 5).Value  = 8
$ws.Cells.Item(115, 6).Value  = "Fruta"
$ws.Cells.Item(115, 7).Value  = 100102
$ws.Cells.Item(115, 8).Value  = "Cítricos"
$ws.Cells.Item(115, 9).Value  = 100102005
$ws.Cells.Item(115, 10).Value = "Naranja"
$ws.Cells.Item(115, 11).Value = "Valencia"
$ws.Cells.Item(115, 12).Value = "Primera"
$ws.Cells.Item(115, 13).Value = 220
$ws.Cells.Item(115, 14).Value = 9000
$ws.Cells.Item(115, 15).Value = 9500
$ws.Cells.Item(115, 16).Value = 9273
$ws.Cells.Item(115, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(115, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(115, 19).Value = 618
$ws.Cells.Item(115, 20).Value = 15
